$wb = $excel.ActiveWorkbook

# --- SQL sheet: remove the trailing test rows (13-15) ---
$wsSql = $wb.Worksheets.Item("SQL")
$wsSql.Rows("13:15").Delete()

# --- Python sheet: remove the trailing test rows (30-32) ---
$wsPython = $wb.Worksheets.Item("Python")
$wsPython.Rows("30:32").Delete()

# --- Links sheet: remove the leading test row (1) ---
$wsLinks = $wb.Worksheets.Item("Links")
$wsLinks.Rows("1:1").Delete()
